$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# --- Insert new columns, preserving row styles (header s=1, data s=2) ---
# 1) New "category" column, inserted right after H (property_category),
#    before the existing I (date) column -> shifts date/legislator_name/legislator_id right by one.
$ws.Columns.Item(9).Insert()

# 2) Two new trailing columns ("source_file", "index") appended after the
#    (now shifted) legislator_id column (L).
$ws.Columns.Item(13).Insert()
$ws.Columns.Item(13).Insert()

# --- Header row (write left-to-right so new shared strings land in the same
#     order the original diff introduces them) ---
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# --- Row 2 (legislator_id 61, 五鼎生技) ---
$ws.Cells.Item(2, 9).Value = "normal"
$ws.Cells.Item(2, 13).Value = "tmpec731"
$ws.Cells.Item(2, 14).Value = 61

# --- Row 3 (legislator_id 62, 萊德科技) ---
$ws.Cells.Item(3, 9).Value = "normal"
$ws.Cells.Item(3, 13).Value = "tmpec731"
$ws.Cells.Item(3, 14).Value = 62
